$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for "2022-Q3", push old "2022-Q2"
#    row down to row 3 (and bump its index value from 0 to 1).
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows.Item(2).Insert()

# Re-apply the numeric-index cell style (bold/bordered/centered) that the
# row-insert did not carry onto the freshly inserted A2 cell; A3 already
# kept it automatically.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.05
$wsTotal.Range("B2").Style = "Normal"
$wsTotal.Range("C2").Style = "Normal"
$wsTotal.Range("D2").Style = "Normal"

$wsTotal.Range("A3").Value = 1

# ---------------------------------------------------------------------
# 2) Quarter sheet: the existing "2022-Q2" tab keeps its current data but
#    is duplicated so the duplicate (placed right after it) becomes the
#    new "2022-Q2" tab, while the original tab is renamed to "2022-Q3"
#    and repopulated with the new quarter's fund holdings.
# ---------------------------------------------------------------------
$wsQ = $wb.Worksheets.Item("2022-Q2")

# Duplicate BEFORE renaming/clearing so the copy preserves the old data.
$wsQ.Copy($null, $wsQ)
$wsDup = $wb.Worksheets.Item("2022-Q2 (2)")

# Free up the "2022-Q2" name on the original sheet, then claim it on the
# duplicate.
$wsQ.Name = "2022-Q3"
$wsDup.Name = "2022-Q2"

# Wipe the old data out of the renamed "2022-Q3" sheet and write the new
# quarter's table.
$wsQ.Cells.Clear()

$wsQ.Range("B1").Value = "基金代码"
$wsQ.Range("C1").Value = "基金名称"
$wsQ.Range("D1").Value = "基金规模"
$wsQ.Range("E1").Value = "股票总仓位"
$wsQ.Range("F1").Value = "仓位占比"
$wsQ.Range("G1").Value = "持有市值(亿元)"
$wsQ.Range("H1").Value = "仓位排名"

$wsQ.Range("A2").Value = 0
$wsQ.Range("B2:G2").NumberFormat = "@"
$wsQ.Range("B2").Value = "166109"
$wsQ.Range("C2").Value = "信澳量化先锋混合（LOF）A"
$wsQ.Range("D2").Value = "0.79"
$wsQ.Range("E2").Value = "88.99"
$wsQ.Range("F2").Value = "4.36"
$wsQ.Range("G2").Value = "0.0344"
$wsQ.Range("H2").Value = 2

$wsQ.Range("A3").Value = 1
$wsQ.Range("B3:G3").NumberFormat = "@"
$wsQ.Range("B3").Value = "003456"
$wsQ.Range("C3").Value = "信澳新目标灵活配置混合"
$wsQ.Range("D3").Value = "0.39"
$wsQ.Range("E3").Value = "94.17"
$wsQ.Range("F3").Value = "2.06"
$wsQ.Range("G3").Value = "0.0080"
$wsQ.Range("H3").Value = 3

$wsQ.Range("A4").Value = 2
$wsQ.Range("B4:G4").NumberFormat = "@"
$wsQ.Range("B4").Value = "166110"
$wsQ.Range("C4").Value = "信澳量化先锋混合（LOF）C"
$wsQ.Range("D4").Value = "0.11"
$wsQ.Range("E4").Value = "88.99"
$wsQ.Range("F4").Value = "4.36"
$wsQ.Range("G4").Value = "0.0048"
$wsQ.Range("H4").Value = 2

# Give the header row + index column the same style used elsewhere in the
# workbook for table headers / index cells, then drop the ad-hoc "@" text
# number-format from the rest of the data cells (keeps values as text
# without leaving a stray style index behind).
$wsTotal.Range("B1").Copy()
$wsQ.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ.Range("A2:A4").PasteSpecial(-4122)
$wsQ.Range("B2:G4").Style = "Normal"

# Keep the original active sheet (workbook-level view metadata is
# untouched by the source diff).
$wsTotal.Activate()

Write-Host "applied 2022-Q3 addition"
